# Rebuilt Panel and loop
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ScintillatorHeight: 145 -> 290
$ws.Range("B8").Value = 290

# MPPCDepth: 0 -> 1
$ws.Range("B16").Value = 1

# New row 24: LoopLength / 140 / mm
$ws.Range("A24").Value = "LoopLength"
$ws.Range("B24").Value = 140
$ws.Range("C24").Value = "mm"

# Move the active selection to A25, matching the saved view state
$ws.Range("A25").Select()
